$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two new rows before row 561, shifting existing rows 561-656 down to 563-658 ---
$ws.Rows("561:562").Insert()

# --- Step 2: populate the two newly-inserted rows (561-562) with what used to be in rows 557-558 ---

# New row 561 <= old row 557 data
$ws.Cells.Item(561, 1).Value = 6
$ws.Cells.Item(561, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(561, 3).Value = "Metropolitana"
$ws.Cells.Item(561, 4).Value = 44481
$ws.Cells.Item(561, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(561, 5).Value = 13
$ws.Cells.Item(561, 6).Value = 100112024
$ws.Cells.Item(561, 7).Value = "Choclo"
$ws.Cells.Item(561, 8).Value = "Choclero"
$ws.Cells.Item(561, 9).Value = "Primera"
$ws.Cells.Item(561, 10).Value = 200
$ws.Cells.Item(561, 11).Value = 45000
$ws.Cells.Item(561, 12).Value = 46000
$ws.Cells.Item(561, 13).Value = 45400
$ws.Cells.Item(561, 14).Value = "$/malla 30 unidades"
$ws.Cells.Item(561, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(561, 16).Value = 1513
$ws.Cells.Item(561, 17).Value = 30
$ws.Cells.Item(561, 18).Value = "Hortaliza"

# New row 562 <= old row 558 data
$ws.Cells.Item(562, 1).Value = 6
$ws.Cells.Item(562, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(562, 3).Value = "Metropolitana"
$ws.Cells.Item(562, 4).Value = 44481
$ws.Cells.Item(562, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(562, 5).Value = 13
$ws.Cells.Item(562, 6).Value = 100112024
$ws.Cells.Item(562, 7).Value = "Choclo"
$ws.Cells.Item(562, 8).Value = "Dulce o Americano"
$ws.Cells.Item(562, 9).Value = "Primera"
$ws.Cells.Item(562, 10).Value = 720
$ws.Cells.Item(562, 11).Value = 44000
$ws.Cells.Item(562, 12).Value = 45000
$ws.Cells.Item(562, 13).Value = 44472
$ws.Cells.Item(562, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(562, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(562, 16).Value = 635
$ws.Cells.Item(562, 17).Value = 70
$ws.Cells.Item(562, 18).Value = "Hortaliza"

# --- Step 3: overwrite rows 557-558 (untouched by the insert) with the new weekly values ---

# Row 557: new date + new price data (H/I/N/O/Q unchanged)
$ws.Cells.Item(557, 4).Value = 44504
$ws.Cells.Item(557, 10).Value = 1300
$ws.Cells.Item(557, 11).Value = 33000
$ws.Cells.Item(557, 12).Value = 35000
$ws.Cells.Item(557, 13).Value = 33769
$ws.Cells.Item(557, 16).Value = 1126

# Row 558: new date + new price data (H/I/N/O/Q unchanged)
$ws.Cells.Item(558, 4).Value = 44504
$ws.Cells.Item(558, 10).Value = 790
$ws.Cells.Item(558, 11).Value = 38000
$ws.Cells.Item(558, 12).Value = 40000
$ws.Cells.Item(558, 13).Value = 39038
$ws.Cells.Item(558, 16).Value = 558
